# Applies the "automatic update" edit: rows 15/16/17 rotate their species-
# observation data (row15<-row16, row16<-row17, row17<-row15 original), and
# rows 22/23 swap their species-observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15 (becomes the old Row 16 observation) ----
$ws.Range("A15").Value = 112381644
$ws.Range("B15").Value = 89557
$ws.Range("E15").Value = 5432
$ws.Range("F15").Value = "Granticka"
$ws.Range("G15").Value = "Porodaedalea chrysoloma"
$ws.Range("H15").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P15").Value = "Väster-Rissjön (Väster-Rissön), Ås lm"
$ws.Range("Q15").Value = 517965
$ws.Range("R15").Value = 7181173
$ws.Range("Z15").Value = "13:20"
$ws.Range("AB15").Value = "13:20"
$ws.Range("AJ15").Value = "gran"
$ws.Range("AK15").Value = "Picea abies"
$ws.Range("AM15").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO15").Value = "Horizontal, dead with ground contact # Picea abies"

# ---- Row 16 (becomes the old Row 17 observation) ----
$ws.Range("A16").Value = 112382121
$ws.Range("B16").Value = 89535
$ws.Range("E16").Value = 1108
$ws.Range("F16").Value = "Harticka"
$ws.Range("G16").Value = "Pelloporus leporinus"
$ws.Range("H16").Value = "(Fr.) Krieglst."
$ws.Range("P16").Value = "Väster-Rissjön (Väster-Rissjön), Ås lm"
$ws.Range("Q16").Value = 517844
$ws.Range("R16").Value = 7181358
$ws.Range("Z16").Value = "13:56"
$ws.Range("AB16").Value = "13:56"

# ---- Row 17 (becomes the original Row 15 observation) ----
$ws.Range("A17").Value = 112375515
$ws.Range("B17").Value = 77637
$ws.Range("E17").Value = 230405
$ws.Range("F17").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G17").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("P17").Value = "Väster-Rissön, Ås lm"
$ws.Range("Q17").Value = 518198
$ws.Range("R17").Value = 7181286
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("AJ17").ClearContents()
$ws.Range("AK17").ClearContents()
$ws.Range("AM17").ClearContents()
$ws.Range("AO17").ClearContents()

# ---- Row 22 (becomes the old Row 23 observation) ----
$ws.Range("A22").Value = 112379312
$ws.Range("B22").Value = 90221
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 3298
$ws.Range("F22").Value = "Trådticka"
$ws.Range("G22").Value = "Climacocystis borealis"
$ws.Range("H22").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Z22").Value = "12:21"
$ws.Range("AB22").Value = "12:21"
$ws.Range("AJ22").Value = "gran"
$ws.Range("AK22").Value = "Picea abies"
$ws.Range("AM22").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO22").Value = "Standing dead tree/snags # Picea abies"

# ---- Row 23 (becomes the original Row 22 observation) ----
$ws.Range("A23").Value = 112375286
$ws.Range("B23").Value = 77637
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 230405
$ws.Range("F23").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G23").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("Z23").ClearContents()
$ws.Range("AB23").ClearContents()
$ws.Range("AJ23").ClearContents()
$ws.Range("AK23").ClearContents()
$ws.Range("AM23").ClearContents()
$ws.Range("AO23").ClearContents()
